$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$players = @(
    @("Eloundou Yvan", "L1C", 698207148),
    @("Kemta Tchoffo", "L1C", 656773516),
    @("Ndoum Junior", "L1C", 653303192),
    @("Bobo Benda Ulrich Gregore", "L1A", 699003156),
    @("Assoumou Assoumou Alahn", "L1A", 696715079)
)

$row = 4
foreach ($player in $players) {
    $ws.Cells.Item($row, 1).Value = $player[0]
    $ws.Cells.Item($row, 2).Value = $player[1]
    $ws.Cells.Item($row, 3).Value = $player[2]
    $row++
}

# Column A grew a bit wider once "Bobo Benda Ulrich Gregore" was added
# (target stored width ~26.14 chars); nudge the column width accordingly.
$ws.Columns.Item(1).ColumnWidth = 25.25

$ws.Range("D8").Select() | Out-Null
